# Auto-generated Excel COM-interop script to refresh market-data columns (H-N)
# on the Odin_Profits workbook, per scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1504.5927
$ws.Range("J19").Value = 1913.0714
$ws.Range("L19").Value = 1913.0714
$ws.Range("N19").Value = -2263.0714
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 955.0833
$ws.Range("I39").Value = 79.111115
$ws.Range("J39").Value = 3583
$ws.Range("K39").Value = 237.333345
$ws.Range("L39").Value = 10749
$ws.Range("M39").Value = 58.66665499999999
$ws.Range("N39").Value = -11341
$ws.Range("H40").Value = 2541.5557
$ws.Range("I40").Value = 2296.75
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 2296.75
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -2121.75
$ws.Range("N40").Value = -4850
$ws.Range("H41").Value = 345.1
$ws.Range("I41").Value = 314.7143
$ws.Range("K41").Value = 314.7143
$ws.Range("M41").Value = 125.2857
$ws.Range("H51").Value = 14661.75
$ws.Range("I51").Value = 7498
$ws.Range("K51").Value = 7498
$ws.Range("M51").Value = -7014
$ws.Range("H80").Value = 617.2941
$ws.Range("I80").Value = 807.375
$ws.Range("J80").Value = 448.33334
$ws.Range("K80").Value = 2422.125
$ws.Range("L80").Value = 1345.00002
$ws.Range("M80").Value = -1424.125
$ws.Range("N80").Value = -3341.00002
$ws.Range("H82").Value = 6547.25
$ws.Range("I82").Value = 506.2
$ws.Range("K82").Value = 1518.6
$ws.Range("M82").Value = -1112.6
$ws.Range("H83").Value = 617.2941
$ws.Range("I83").Value = 807.375
$ws.Range("J83").Value = 448.33334
$ws.Range("K83").Value = 7266.375
$ws.Range("L83").Value = 4035.00006
$ws.Range("M83").Value = -2274.375
$ws.Range("N83").Value = -14019.00006
$ws.Range("H85").Value = 6547.25
$ws.Range("I85").Value = 506.2
$ws.Range("K85").Value = 1518.6
$ws.Range("M85").Value = -114.5999999999999
$ws.Range("H98").Value = 4321.32
$ws.Range("I98").Value = 4122
$ws.Range("K98").Value = 4122
$ws.Range("M98").Value = -2624
$ws.Range("H113").Value = 3216.318
$ws.Range("J113").Value = 3299.5625
$ws.Range("L113").Value = 3299.5625
$ws.Range("N113").Value = -9807.5625
$ws.Range("H122").Value = 4321.32
$ws.Range("I122").Value = 4122
$ws.Range("K122").Value = 12366
$ws.Range("M122").Value = -9916
$ws.Range("H132").Value = 521891.56
$ws.Range("I132").Value = 681541.5
$ws.Range("K132").Value = 2044624.5
$ws.Range("M132").Value = -2042094.5
$ws.Range("H138").Value = 3722.0588
$ws.Range("J138").Value = 4339.2974
$ws.Range("L138").Value = 13017.8922
$ws.Range("N138").Value = -23297.8922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6122.7144
$ws.Range("I2").Value = 3649.5217
$ws.Range("K2").Value = 3649.5217
$ws.Range("M2").Value = -3536.5217
$ws.Range("H32").Value = 1588506.5
$ws.Range("I32").Value = 876.5806
$ws.Range("K32").Value = 876.5806
$ws.Range("M32").Value = -589.5806
$ws.Range("H45").Value = 1837.9688
$ws.Range("I45").Value = 1678.8518
$ws.Range("K45").Value = 1678.8518
$ws.Range("M45").Value = -1301.8518
$ws.Range("H61").Value = 5422.757
$ws.Range("I61").Value = 8331
$ws.Range("K61").Value = 8331
$ws.Range("M61").Value = -8119
$ws.Range("H116").Value = 6122.7144
$ws.Range("I116").Value = 3649.5217
$ws.Range("K116").Value = 3649.5217
$ws.Range("M116").Value = -1355.5217
$ws.Range("H122").Value = 2964.4043
$ws.Range("I122").Value = 2512.1667
$ws.Range("K122").Value = 7536.500100000001
$ws.Range("M122").Value = -5086.500100000001
$ws.Range("H136").Value = 5422.757
$ws.Range("I136").Value = 8331
$ws.Range("K136").Value = 24993
$ws.Range("M136").Value = -22443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6122.7144
$ws.Range("I3").Value = 3649.5217
$ws.Range("K3").Value = 3649.5217
$ws.Range("M3").Value = -3535.5217
$ws.Range("H94").Value = 5575.6816
$ws.Range("I94").Value = 980.9167
$ws.Range("K94").Value = 980.9167
$ws.Range("M94").Value = -529.9167
$ws.Range("H99").Value = 4979.478
$ws.Range("I99").Value = 3802.1333
$ws.Range("J99").Value = 7187
$ws.Range("K99").Value = 3802.1333
$ws.Range("L99").Value = 7187
$ws.Range("M99").Value = -2304.1333
$ws.Range("N99").Value = -10183
$ws.Range("H134").Value = 1332756.4
$ws.Range("I134").Value = 1543436.5
$ws.Range("J134").Value = 26540.2
$ws.Range("K134").Value = 4630309.5
$ws.Range("L134").Value = 79620.60000000001
$ws.Range("M134").Value = -4627774.5
$ws.Range("N134").Value = -84690.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3922.3777
$ws.Range("I31").Value = 1272.3158
$ws.Range("J31").Value = 5858.9614
$ws.Range("K31").Value = 1272.3158
$ws.Range("L31").Value = 5858.9614
$ws.Range("M31").Value = -977.3158000000001
$ws.Range("N31").Value = -6448.9614
$ws.Range("H34").Value = 3922.3777
$ws.Range("I34").Value = 1272.3158
$ws.Range("J34").Value = 5858.9614
$ws.Range("K34").Value = 1272.3158
$ws.Range("L34").Value = 5858.9614
$ws.Range("M34").Value = -1070.3158
$ws.Range("N34").Value = -6262.9614
$ws.Range("H58").Value = 43487264
$ws.Range("I58").Value = 111120040
$ws.Range("J58").Value = 9051.286
$ws.Range("K58").Value = 111120040
$ws.Range("L58").Value = 9051.286
$ws.Range("M58").Value = -111119837
$ws.Range("N58").Value = -9457.286
$ws.Range("H62").Value = 4727.609
$ws.Range("I62").Value = 2486.1667
$ws.Range("J62").Value = 5518.706
$ws.Range("K62").Value = 2486.1667
$ws.Range("L62").Value = 5518.706
$ws.Range("M62").Value = -1862.1667
$ws.Range("N62").Value = -6766.706
$ws.Range("H65").Value = 4727.609
$ws.Range("I65").Value = 2486.1667
$ws.Range("J65").Value = 5518.706
$ws.Range("K65").Value = 12430.8335
$ws.Range("L65").Value = 27593.53
$ws.Range("M65").Value = -9310.833500000001
$ws.Range("N65").Value = -33833.53
$ws.Range("H107").Value = 520.2353000000001
$ws.Range("I107").Value = 521.5
$ws.Range("K107").Value = 521.5
$ws.Range("M107").Value = 1398.5
$ws.Range("H134").Value = 22229212
$ws.Range("I134").Value = 25646798
$ws.Range("K134").Value = 76940394
$ws.Range("M134").Value = -76937859
$ws.Range("H136").Value = 43487264
$ws.Range("I136").Value = 111120040
$ws.Range("J136").Value = 9051.286
$ws.Range("K136").Value = 333360120
$ws.Range("L136").Value = 27153.858
$ws.Range("M136").Value = -333357570
$ws.Range("N136").Value = -32253.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 512.8570999999999
$ws.Range("I7").Value = 548.3333
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1644.9999
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -1532.9999
$ws.Range("N7").Value = -1124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 99497
$ws.Range("J119").Value = 99497
$ws.Range("L119").Value = 99497
$ws.Range("N119").Value = -109173
$ws.Range("H121").Value = 87333.336
$ws.Range("J121").Value = 87333.336
$ws.Range("L121").Value = 87333.336
$ws.Range("N121").Value = -90827.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8317.885
$ws.Range("I7").Value = 6790.3335
$ws.Range("J7").Value = 10400.909
$ws.Range("K7").Value = 6790.3335
$ws.Range("L7").Value = 10400.909
$ws.Range("M7").Value = -6678.3335
$ws.Range("N7").Value = -10624.909
$ws.Range("H40").Value = 3944.4
$ws.Range("I40").Value = 2840.2666
$ws.Range("J40").Value = 7256.8
$ws.Range("K40").Value = 2840.2666
$ws.Range("L40").Value = 7256.8
$ws.Range("M40").Value = -2704.2666
$ws.Range("N40").Value = -7528.8
$ws.Range("H61").Value = 8005.289
$ws.Range("J61").Value = 10168.733
$ws.Range("L61").Value = 10168.733
$ws.Range("N61").Value = -10572.733
$ws.Range("H113").Value = 8005.289
$ws.Range("J113").Value = 10168.733
$ws.Range("L113").Value = 10168.733
$ws.Range("N113").Value = -14508.733
$ws.Range("H122").Value = 5584
$ws.Range("I122").Value = 3817.7
$ws.Range("K122").Value = 11453.1
$ws.Range("M122").Value = -9003.099999999999
$ws.Range("H126").Value = 8317.885
$ws.Range("I126").Value = 6790.3335
$ws.Range("J126").Value = 10400.909
$ws.Range("K126").Value = 20371.0005
$ws.Range("L126").Value = 31202.727
$ws.Range("M126").Value = -17901.0005
$ws.Range("N126").Value = -36142.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22624.75
$ws.Range("I62").Value = 20375
$ws.Range("K62").Value = 20375
$ws.Range("M62").Value = -19751
$ws.Range("H65").Value = 22624.75
$ws.Range("I65").Value = 20375
$ws.Range("K65").Value = 101875
$ws.Range("M65").Value = -98755
$ws.Range("H100").Value = 1307.625
$ws.Range("I100").Value = 1422.2
$ws.Range("K100").Value = 2844.4
$ws.Range("M100").Value = -2303.4
$ws.Range("H107").Value = 822.0323
$ws.Range("I107").Value = 557.55
$ws.Range("J107").Value = 1302.909
$ws.Range("K107").Value = 1672.65
$ws.Range("L107").Value = 3908.727
$ws.Range("M107").Value = 247.3500000000001
$ws.Range("N107").Value = -7748.727000000001
$ws.Range("H126").Value = 4931.2383
$ws.Range("I126").Value = 4928.4546
$ws.Range("K126").Value = 14785.3638
$ws.Range("M126").Value = -12315.3638
$ws.Range("H136").Value = 11909556
$ws.Range("I136").Value = 15628694
$ws.Range("K136").Value = 46886082
$ws.Range("M136").Value = -46883532

Write-Host "Odin_Profits market data refresh applied."